$wb = $excel.ActiveWorkbook

# --- Sheet "Subs": countUp tasks should not have a time value ---
$wsSubs = $wb.Worksheets.Item("Subs")
$wsSubs.Range("C5").Value = 2
$wsSubs.Range("E5").Value = 0
[void]$wsSubs.Range("F5").Select()

# --- Sheet "Subsubs": re-split the Quick Det rows, add Research + Teaching subtasks ---
$wsSubsubs = $wb.Worksheets.Item("Subsubs")

# Row 5: Quick Det / LN Thm 3.5 (was row 4 content)
$wsSubsubs.Range("A5").Value = "Quick Det"
$wsSubsubs.Range("B5").Value = "LN Thm 3.5"
$wsSubsubs.Range("C5").Value = 2

# Row 6: Quick Det / Lm-Thm 3.7-8 (was row 5 content)
$wsSubsubs.Range("A6").Value = "Quick Det"
$wsSubsubs.Range("B6").Value = "Lm-Thm 3.7-8"
$wsSubsubs.Range("C6").Value = 4

# Row 7: new Teaching subtask "Cours S5"
$wsSubsubs.Range("A7").Value = "Teaching"
$wsSubsubs.Range("B7").Value = "Cours S5"
$wsSubsubs.Range("C7").Value = 2.5

# Row 8: new Teaching subtask "TP S5"
$wsSubsubs.Range("A8").Value = "Teaching"
$wsSubsubs.Range("B8").Value = "TP S5"
$wsSubsubs.Range("C8").Value = 2.5

# --- Sheet "Tasks": update the daily time stamp description and move selection ---
$wsTasks = $wb.Worksheets.Item("Tasks")
$wsTasks.Range("C2").Value = "Daily time stamp for 24/10/2021"
[void]$wsTasks.Range("C2").Select()

# Row 4 becomes a Research subtask ("Distr convolution")
$wsSubsubs.Range("A4").Value = "Research"
$wsSubsubs.Range("B4").Value = "Distr convolution"
$wsSubsubs.Range("C4").Value = 3.5

[void]$wsSubsubs.Activate()
